$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "marks" column header
$ws.Range("E1").Value = "marks"

# Update existing "answer" scores in column D
$ws.Range("D2").Value = 6
$ws.Range("D3").Value = 5
$ws.Range("D5").Value = 5
$ws.Range("D6").Value = 6

# Fill in marks column for the new rows
$ws.Range("E2").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("E4").Value = 0

# Update the selection to match the recorded view state
$ws.Range("F5").Select()
